$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「数学者アル＝フワーリズミー」" (row 57) was removed from the
# spreadsheet. Deleting the entire row shifts every subsequent row up by
# one, which also updates the sheet's used-range dimension from
# A1:C232 to A1:C231 automatically.
$ws.Rows("57").Delete()
